$d = $word.ActiveDocument

# --- Paragraph "1) The percentage of male players..." : merge its several
#     runs into a single run (text itself is unchanged). ---
$p2 = $d.Paragraphs.Item(2)
$full2 = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$full2.Delete()
$ins2 = $d.Range($p2.Range.Start, $p2.Range.Start)
$ins2.InsertAfter("1) The percentage of male players in the game is 84%, which is very high as compared to the female player percentage, which is 14%. Looks like the game is really famous with the male players.")

# --- Paragraph "2) The age group most interested..." : merge its several
#     runs into a single run (text itself is unchanged). ---
$p4 = $d.Paragraphs.Item(4)
$full4 = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$full4.Delete()
$ins4 = $d.Range($p4.Range.Start, $p4.Range.Start)
$ins4.InsertAfter("2) The age group most interested in the game is between 20-24years, so if the company owning the game would want to increase sales, it should target this age group.")

# --- Paragraph "3)Although the most popular item..." : split into three
#     paragraphs -- the original sentence (as one merged run), a new
#     "4) The top spender ..." paragraph, and a trailing empty paragraph
#     that keeps the _GoBack bookmark. ---

# Step A: split the paragraph right at the (hidden) _GoBack bookmark so the
# bookmark, together with the little bit of trailing text it wraps, becomes
# its own paragraph.
$bm = $d.Bookmarks("_GoBack")
$bmStart = $bm.Start
$r = $d.Range($bmStart, $bmStart)
$r.InsertParagraphAfter()

# Step B: move the trailing text ("tems too.") out of the new
# bookmark-paragraph and back onto the end of the sentence paragraph, so the
# bookmark paragraph ends up containing only the bookmark.
$p6 = $d.Paragraphs.Item(6)
$p7 = $d.Paragraphs.Item(7)
$bm2 = $d.Bookmarks("_GoBack")
$afterBookmarkRange = $d.Range($bm2.End, $p7.Range.End - 1)
$afterBookmarkRange.Delete()
$p6 = $d.Paragraphs.Item(6)
$insertPoint = $d.Range($p6.Range.End - 1, $p6.Range.End - 1)
$insertPoint.InsertAfter("tems too.")

# Step C: merge the sentence paragraph's many runs into a single run.
$p6 = $d.Paragraphs.Item(6)
$full6 = $d.Range($p6.Range.Start, $p6.Range.End - 1)
$full6.Delete()
$ins6 = $d.Range($p6.Range.Start, $p6.Range.Start)
$ins6.InsertAfter("3)Although the most popular item is Final Critic, there is not much difference in the purchase count of Final Critic with the second most popular item, Oathbreaker, last Hope of the Breaking Storm. However, the total purchase value differs by almost `$9. Company should try to promote both items with great efforts as these are the most profitable items too.")

# Step D: insert the new "4) The top spender ..." paragraph right after the
# sentence paragraph (i.e. before the now bookmark-only paragraph).
$p6 = $d.Paragraphs.Item(6)
$p6.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs.Item(7)
$ins7 = $d.Range($p7.Range.Start, $p7.Range.Start)
$ins7.InsertAfter("4) The top spender is Lisosia93 with 5 purchase count and total purchase value of `$18.96.")
